$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row to uppercase labels
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "IP"
$ws.Range("D1").Value = "Port"
$ws.Range("E1").Value = "Status"

# Reorder row 2 data to match the ID, Name, IP, Port, Status column layout
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "My_Viettel"
$ws.Range("C2").Value = "192.168.1.1"
$ws.Range("D2").Value = 3000
$ws.Range("E2").Value = $false
